$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.46 = 30277.02 pesos`n✅ 30277.02 pesos = 7.39 = 964.21 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet: rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 133.996
$wsTasas.Range("O10").Value = 4057
$wsTasas.Range("N12").Value = 4097.8
$wsTasas.Range("O12").Value = 130.5
